# ============================================================================
# feat: add 2022-Q1 data
#
# The workbook's old "总计" (aggregate/summary) sheet is renamed to "2022-Q1"
# and repopulated with the quarter's fund-holdings table (same A:H layout as the
# other quarterly sheets). A fresh "总计" sheet is inserted right after it, carrying
# the same date/count/value summary as before plus a new leading row for 2022-Q1.
# ============================================================================

$wb = $excel.ActiveWorkbook

# --- Step 1: the current "总计" sheet (sheetId 6) becomes "2022-Q1". -------------
# We keep the sheet object (so sheetId/position stay put) and overwrite its data.
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# --- Step 2: insert a brand-new "总计" sheet right after "2022-Q1". --------------
# Copying "2021-Q4" gives the new sheet the same header/column styling used by
# every quarterly sheet; we'll clear it down to the small A:D summary shape after.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $q1)
$total = $wb.Worksheets.Item("2021-Q4 (2)")
$total.Name = "总计"

# ============================================================================
# Populate "2022-Q1" with the fund-holdings table (header + 29 funds).
# ============================================================================
$q1.Range("A1:H31").Clear()

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"
$q1.Range("B1:H1").Font.Bold = $true

# Text-format the columns that must stay plain text (codes / decimals with
# trailing zeros that a numeric cell would otherwise drop), matching the
# source data's inlineStr cell types.
$q1.Range("B2:G30").NumberFormat = "@"

$q1FundData = @(
    ,@("000729", "建信中小盘先锋股票", "33.97", "89.28", "4.56", "1.5490", 6)
    ,@("070021", "嘉实主题新动力混合", "24.04", "93.93", "4.68", "1.1251", 8)
    ,@("530005", "建信优化配置混合", "21.57", "86.97", "4.54", "0.9793", 4)
    ,@("011637", "广发沪港深价值成长混合型证券投资基金A", "12.44", "92.96", "6.04", "0.7514", 5)
    ,@("010452", "广发瑞福精选混合A", "16.29", "78.69", "4.10", "0.6679", 3)
    ,@("000985", "嘉实逆向策略股票", "13.64", "93.90", "4.65", "0.6343", 8)
    ,@("010330", "东吴兴享成长混合A", "11.63", "80.15", "5.38", "0.6257", 3)
    ,@("001070", "建信信息产业股票", "10.67", "81.80", "4.54", "0.4844", 4)
    ,@("000756", "建信潜力新蓝筹股票", "10.47", "84.61", "4.53", "0.4743", 5)
    ,@("005598", "广发中小盘精选混合", "9.64", "93.99", "4.89", "0.4714", 4)
    ,@("001305", "九泰天富改革新动力混合A", "3.74", "88.86", "8.53", "0.3190", 2)
    ,@("001782", "九泰久益灵活配置混合A", "2.33", "94.33", "9.48", "0.2209", 2)
    ,@("008962", "建信科技创新混合A", "4.34", "84.82", "4.28", "0.1858", 5)
    ,@("001844", "九泰久益灵活配置混合C", "1.47", "94.33", "9.48", "0.1394", 2)
    ,@("000308", "建信创新中国混合", "3.11", "84.50", "3.99", "0.1241", 5)
    ,@("000717", "融通转型三动力灵活配置混合A", "3.83", "94.89", "3.20", "0.1226", 10)
    ,@("010453", "广发瑞福精选混合C", "2.42", "78.69", "4.10", "0.0992", 3)
    ,@("000969", "前海开源大安全核心精选灵活配置混合", "1.39", "91.04", "4.29", "0.0596", 4)
    ,@("009912", "九泰天富改革新动力混合C", "0.59", "88.86", "8.53", "0.0503", 2)
    ,@("001060", "前海开源高端装备制造灵活配置混合", "0.97", "89.88", "4.29", "0.0416", 4)
    ,@("004332", "恒生前海沪港深新兴产业精选混合", "0.52", "80.98", "4.57", "0.0238", 2)
    ,@("001734", "广发百发大数据策略成长灵活配置混合A", "0.65", "85.72", "3.21", "0.0209", 3)
    ,@("001735", "广发百发大数据策略成长灵活配置混合E", "0.65", "85.72", "3.21", "0.0209", 3)
    ,@("011638", "广发沪港深价值成长混合型证券投资基金C", "0.33", "92.96", "6.04", "0.0199", 5)
    ,@("009828", "融通转型三动力灵活配置混合C", "0.59", "94.89", "3.20", "0.0189", 10)
    ,@("011462", "东吴兴享成长混合C", "0.33", "80.15", "5.38", "0.0178", 3)
    ,@("008963", "建信科技创新混合C", "0.26", "84.82", "4.28", "0.0111", 5)
    ,@("008437", "九泰行业优选灵活配置混合A", "0.11", "51.13", "7.15", "0.0079", 1)
    ,@("008438", "九泰行业优选灵活配置混合C", "0.06", "51.13", "7.15", "0.0043", 1)
)

for ($i = 0; $i -lt $q1FundData.Count; $i++) {
    $row = $q1FundData[$i]
    $r = $i + 2
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# Column A (row index) keeps the bold/bordered header style used throughout
# the workbook's quarterly sheets.
$q1.Range("A2:A30").Font.Bold = $true

# ============================================================================
# Populate "总计" with the date/count/value summary (header + 6 rows).
# ============================================================================
$total.Range("A1:H31").Clear()

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$total.Range("B1:D1").Font.Bold = $true

$totalData = @(
    ,@("2022-Q1", 29, 9.27)
    ,@("2021-Q4", 30, 14)
    ,@("2021-Q3", 12, 2.94)
    ,@("2021-Q2", 5, 1.27)
    ,@("2021-Q1", 13, 4.64)
    ,@("2020-Q4", 8, 6.5)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $row = $totalData[$i]
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

$total.Range("A2:A7").Font.Bold = $true

$wb.Worksheets.Item(1).Activate()

